$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-259) holds the "Förändrad" (last changed) date as a
# date-serial number. Every row's value moves from 45172 to 45175
# (i.e. 2023-09-03 -> 2023-09-06), leaving formatting/style untouched.
$ws.Range("C2:C259").Value = 45175
